# Auto-generated edit script applying the Ixion_Profits.xlsx diff
# Updates cached numeric values (no formulas in this workbook) across 8 sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 12410
$ws.Range("I101").Value = 250
$ws.Range("J101").Value = 15450
$ws.Range("K101").Value = 750
$ws.Range("L101").Value = 46350
$ws.Range("M101").Value = 872
$ws.Range("N101").Value = -49594
$ws.Range("H129").Value = 902.2782999999999
$ws.Range("I129").Value = 507
$ws.Range("J129").Value = 947.7126500000001
$ws.Range("K129").Value = 1521
$ws.Range("L129").Value = 2843.13795
$ws.Range("M129").Value = 3479
$ws.Range("N129").Value = -12843.13795
$ws.Range("H132").Value = 1895.0869
$ws.Range("I132").Value = 1694.3334
$ws.Range("K132").Value = 5083.0002
$ws.Range("M132").Value = -2553.0002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 140.25
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -432
$ws.Range("H32").Value = 5476.75
$ws.Range("I32").Value = 4430.9
$ws.Range("J32").Value = 8962.916999999999
$ws.Range("K32").Value = 4430.9
$ws.Range("L32").Value = 8962.916999999999
$ws.Range("M32").Value = -4143.9
$ws.Range("N32").Value = -9536.916999999999
$ws.Range("H132").Value = 3692.2646
$ws.Range("I132").Value = 1898.5264
$ws.Range("J132").Value = 5964.3335
$ws.Range("K132").Value = 5695.5792
$ws.Range("L132").Value = 17893.0005
$ws.Range("M132").Value = -3165.5792
$ws.Range("N132").Value = -22953.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4997.5835
$ws.Range("I134").Value = 6853.7144
$ws.Range("K134").Value = 20561.1432
$ws.Range("M134").Value = -18026.1432

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2081.7144
$ws.Range("I31").Value = 1004.44446
$ws.Range("J31").Value = 4020.8
$ws.Range("K31").Value = 1004.44446
$ws.Range("L31").Value = 4020.8
$ws.Range("M31").Value = -709.44446
$ws.Range("N31").Value = -4610.8
$ws.Range("H34").Value = 2081.7144
$ws.Range("I34").Value = 1004.44446
$ws.Range("J34").Value = 4020.8
$ws.Range("K34").Value = 1004.44446
$ws.Range("L34").Value = 4020.8
$ws.Range("M34").Value = -802.44446
$ws.Range("N34").Value = -4424.8
$ws.Range("H86").Value = 2114.2222
$ws.Range("I86").Value = 2197.9285
$ws.Range("J86").Value = 1821.25
$ws.Range("K86").Value = 2197.9285
$ws.Range("L86").Value = 1821.25
$ws.Range("M86").Value = -1074.9285
$ws.Range("N86").Value = -4067.25
$ws.Range("H89").Value = 2114.2222
$ws.Range("I89").Value = 2197.9285
$ws.Range("J89").Value = 1821.25
$ws.Range("K89").Value = 10989.6425
$ws.Range("L89").Value = 9106.25
$ws.Range("M89").Value = -5373.6425
$ws.Range("N89").Value = -20338.25
$ws.Range("H139").Value = 40453.332
$ws.Range("J139").Value = 40453.332
$ws.Range("L139").Value = 40453.332
$ws.Range("N139").Value = -50733.332

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 400395.28
$ws.Range("I5").Value = 616.125
$ws.Range("J5").Value = 857285.7
$ws.Range("K5").Value = 1848.375
$ws.Range("L5").Value = 2571857.1
$ws.Range("M5").Value = -1736.375
$ws.Range("N5").Value = -2572081.1
$ws.Range("H55").Value = 3137.375
$ws.Range("J55").Value = 3514.1428
$ws.Range("L55").Value = 10542.4284
$ws.Range("N55").Value = -10896.4284
$ws.Range("H62").Value = 6040
$ws.Range("J62").Value = 6040
$ws.Range("L62").Value = 18120
$ws.Range("N62").Value = -19492
$ws.Range("H63").Value = 5725
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 5725
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 17175
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -18673
$ws.Range("H64").Value = 3085.1667
$ws.Range("I64").Value = 1755.5
$ws.Range("J64").Value = 3750
$ws.Range("K64").Value = 5266.5
$ws.Range("L64").Value = 11250
$ws.Range("M64").Value = -4996.5
$ws.Range("N64").Value = -11790
$ws.Range("H65").Value = 6040
$ws.Range("J65").Value = 6040
$ws.Range("L65").Value = 54360
$ws.Range("N65").Value = -61224
$ws.Range("H66").Value = 5725
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 5725
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 51525
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -59013
$ws.Range("H67").Value = 3085.1667
$ws.Range("I67").Value = 1755.5
$ws.Range("J67").Value = 3750
$ws.Range("K67").Value = 5266.5
$ws.Range("L67").Value = 11250
$ws.Range("M67").Value = -4330.5
$ws.Range("N67").Value = -13122
$ws.Range("H68").Value = 2015.7142
$ws.Range("I68").Value = 477.5
$ws.Range("J68").Value = 4066.6667
$ws.Range("K68").Value = 1432.5
$ws.Range("L68").Value = 12200.0001
$ws.Range("M68").Value = -621.5
$ws.Range("N68").Value = -13822.0001
$ws.Range("H71").Value = 2015.7142
$ws.Range("I71").Value = 477.5
$ws.Range("J71").Value = 4066.6667
$ws.Range("K71").Value = 4297.5
$ws.Range("L71").Value = 36600.0003
$ws.Range("M71").Value = -241.5
$ws.Range("N71").Value = -44712.0003
$ws.Range("H113").Value = 333837.12
$ws.Range("I113").Value = 487.65
$ws.Range("J113").Value = 1000536.1
$ws.Range("K113").Value = 1462.95
$ws.Range("L113").Value = 3001608.3
$ws.Range("M113").Value = 707.0500000000002
$ws.Range("N113").Value = -3005948.3
$ws.Range("H131").Value = 1516105.2
$ws.Range("J131").Value = 978.2105
$ws.Range("L131").Value = 2934.6315
$ws.Range("N131").Value = -13014.6315
$ws.Range("H134").Value = 9959.022000000001
$ws.Range("I134").Value = 9933.083000000001
$ws.Range("J134").Value = 9968.75
$ws.Range("K134").Value = 29799.249
$ws.Range("L134").Value = 29906.25
$ws.Range("M134").Value = -24729.249
$ws.Range("N134").Value = -40046.25
$ws.Range("H135").Value = 400395.28
$ws.Range("I135").Value = 616.125
$ws.Range("J135").Value = 857285.7
$ws.Range("K135").Value = 5545.125
$ws.Range("L135").Value = 7715571.3
$ws.Range("M135").Value = -3010.125
$ws.Range("N135").Value = -7720641.3
$ws.Range("H138").Value = 11860
$ws.Range("I138").Value = 11860
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 35580
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -30440
$ws.Range("N138").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H130").Value = 47128.57
$ws.Range("J130").Value = 47128.57
$ws.Range("L130").Value = 47128.57
$ws.Range("N130").Value = -57168.57

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 54736.42
$ws.Range("I7").Value = 73213.71000000001
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 73213.71000000001
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -73101.71000000001
$ws.Range("N7").Value = -3224
$ws.Range("H82").Value = 80213.84
$ws.Range("I82").Value = 1280
$ws.Range("K82").Value = 1280
$ws.Range("M82").Value = -919
$ws.Range("H85").Value = 80213.84
$ws.Range("I85").Value = 1280
$ws.Range("K85").Value = 1280
$ws.Range("M85").Value = -32
$ws.Range("H114").Value = 43000
$ws.Range("J114").Value = 43000
$ws.Range("L114").Value = 43000
$ws.Range("N114").Value = -51678
$ws.Range("H115").Value = 45000
$ws.Range("J115").Value = 45000
$ws.Range("L115").Value = 45000
$ws.Range("N115").Value = -47350
$ws.Range("H126").Value = 54736.42
$ws.Range("I126").Value = 73213.71000000001
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 219641.13
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -217171.13
$ws.Range("N126").Value = -13940

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 49084.668
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H35").Value = 49084.668
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H81").Value = 2261.25
$ws.Range("I81").Value = 1348.3334
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 2696.6668
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -1635.6668
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 2261.25
$ws.Range("I84").Value = 1348.3334
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 13483.334
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -8179.333999999999
$ws.Range("N84").Value = -60608
$ws.Range("H135").Value = 47500
$ws.Range("J135").Value = 47500
$ws.Range("L135").Value = 47500
$ws.Range("N135").Value = -57640
